$d = $word.ActiveDocument

# Update the date heading (first paragraph, centered, Arial 30)
$d.Paragraphs.Item(1).Range.Text = "2025-08-29 Friday"

# Update the answer cells in the single results table, addressed by
# (row, column) to avoid any ambiguity from repeated/overlapping values.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "75÷6=12, 3"
$t.Cell(1, 2).Range.Text = "66÷5=13, 1"
$t.Cell(1, 3).Range.Text = "12÷2=6, 0"
$t.Cell(1, 4).Range.Text = "94÷6=15, 4"
$t.Cell(1, 5).Range.Text = "82÷3=27, 1"
$t.Cell(5, 1).Range.Text = "73÷7=10, 3"
$t.Cell(5, 2).Range.Text = "52÷8=6, 4"
$t.Cell(5, 3).Range.Text = "87÷3=29, 0"
$t.Cell(5, 4).Range.Text = "55÷9=6, 1"
$t.Cell(5, 5).Range.Text = "77÷5=15, 2"
$t.Cell(9, 1).Range.Text = "82÷5=16, 2"
$t.Cell(9, 2).Range.Text = "49÷7=7, 0"
$t.Cell(9, 3).Range.Text = "61÷5=12, 1"
$t.Cell(9, 4).Range.Text = "97÷9=10, 7"
$t.Cell(9, 5).Range.Text = "50÷4=12, 2"
$t.Cell(13, 1).Range.Text = "88÷3=29, 1"
$t.Cell(13, 2).Range.Text = "18÷4=4, 2"
$t.Cell(13, 3).Range.Text = "34÷8=4, 2"
$t.Cell(13, 4).Range.Text = "39÷5=7, 4"
$t.Cell(13, 5).Range.Text = "12÷9=1, 3"
$t.Cell(17, 1).Range.Text = "76÷8=9, 4"
$t.Cell(17, 2).Range.Text = "72÷6=12, 0"
$t.Cell(17, 3).Range.Text = "17÷7=2, 3"
$t.Cell(17, 4).Range.Text = "19÷6=3, 1"
$t.Cell(17, 5).Range.Text = "83÷6=13, 5"
